# Applies the commit "update appointment and appointmentResponse" changes:
#  - Metadata sheet: Version 0.1.0 -> 0.2.0, Status active -> draft,
#    Date 2023-04-27T14:52:10+02:00 -> 2023-07-20T11:38:03+02:00
#  - Rename "Include from virtual-service-" sheet -> "Include from unknown"
#    and update its "System URI" value to the new canonical OID
#  - Rename "Include from contact-point-sy" sheet -> "Include from ContactPointSyst"
#    and update its "System URI" value to the new canonical URL
#  - Restore wrap-text/top-vertical alignment formatting on the header and
#    body rows of every sheet (re-asserts applyAlignment on those styles)

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -----------------------------------------------------
$meta = $wb.Worksheets.Item(1)
$meta.Cells.Item(3, 2).Value = "0.2.0"
$meta.Cells.Item(6, 2).Value = "draft"
$meta.Cells.Item(8, 2).Value = "2023-07-20T11:38:03+02:00"

# --- Include from virtual-service- (2nd sheet) ---------------------------
$vs = $wb.Worksheets.Item(2)
$vs.Cells.Item(4, 2).Value = "urn:oid:2.16.840.1.113883.4.642.4.1809"
$vs.Name = "Include from unknown"

# --- Include from contact-point-sy (3rd sheet) ---------------------------
$cp = $wb.Worksheets.Item(3)
$cp.Cells.Item(4, 2).Value = "http://hl7.org/fhir/contact-point-system"
$cp.Name = "Include from ContactPointSyst"

# --- Re-assert wrap-text / top alignment on the existing cells -----------
# (Keeps the bold/filled header style and the bordered body style intact,
#  but re-applies the alignment so it is flagged "applied" again.)
$meta.Range("A1:B1").WrapText = $true
$meta.Range("A2:B14").WrapText = $true

$vs.Range("A1").WrapText = $true
$vs.Range("A2").WrapText = $true
$vs.Range("A3:B3").WrapText = $true
$vs.Range("A4:B4").WrapText = $true

$cp.Range("A1").WrapText = $true
$cp.Range("A2").WrapText = $true
$cp.Range("A3:B3").WrapText = $true
$cp.Range("A4:B4").WrapText = $true
